# Apply updated crypto price / 1h-volume data to Sheet1 (rows 2-51).
# Column D = Price, Column E = Volume(1h) percentage change.
# NumberFormat is forced to "@" (text) before assignment so that
# numeric-looking strings (e.g. "5.21", "36.80") are stored as
# literal text, matching the source data's formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.205.74'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.23%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.517.97'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.89%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '519.35'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '131.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.555'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -0.43%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.517.92'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0974'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.86%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.21'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.35%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.21%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.935.95'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.189.87'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '22.19'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.57%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.504.43'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.73'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.34%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '323.74'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.01%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.08%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.03'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +4.92%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.996'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.42%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.84'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.38%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.406'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.49%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +1.31%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -1.47%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.37'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.32%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0749'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '167.47'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.01%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.43%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.27'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -0.47%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.16%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.998'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '18.04'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.12%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.27'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.54%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.89%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.80'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -1.44%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.778'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.41%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '277.31'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +1.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.10'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +1.82%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.599'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0921'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '122.04'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.48%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0499'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '17.74'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.64%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.11%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '17.00'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.45%  '
